$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 38:43 ("Atividade" column) were stored as text "6" before; the
# refreshed export writes them as a genuine number 6 (display unchanged).
for ($r = 38; $r -le 43; $r++) {
    $ws.Range("A$r").Value = 6
}

# New rows 44:49 - SVC / "Branch and Bound - Sem Normalização" results
# (Atividade stored as text "7", matching the source export's dtype).
$ws.Range("A44:A49").NumberFormat = "@"

$ws.Range("A44").Value = "7"
$ws.Range("B44").Value = "SVC"
$ws.Range("C44").Value = "Branch and Bound - Sem Normalização - Desbalanceado"
$ws.Range("D44").Value = "('rbf', 10, '0.9201')"
$ws.Range("E44").Value = "Accuracy"
$ws.Range("F44").Value = 0.9382352941176469

$ws.Range("A45").Value = "7"
$ws.Range("B45").Value = "SVC"
$ws.Range("C45").Value = "Branch and Bound - Sem Normalização - Desbalanceado"
$ws.Range("D45").Value = "('rbf', 10, '0.9201')"
$ws.Range("E45").Value = "Recall"
$ws.Range("F45").Value = 0.08333333333333333

$ws.Range("A46").Value = "7"
$ws.Range("B46").Value = "SVC"
$ws.Range("C46").Value = "Branch and Bound - Sem Normalização - Desbalanceado"
$ws.Range("D46").Value = "('rbf', 10, '0.9201')"
$ws.Range("E46").Value = "Precision"
$ws.Range("F46").Value = 0.35

$ws.Range("A47").Value = "7"
$ws.Range("B47").Value = "SVC"
$ws.Range("C47").Value = "Branch and Bound - Sem Normalização - Balanceado"
$ws.Range("D47").Value = "('sigmoid', 1, '0.0201')"
$ws.Range("E47").Value = "Accuracy"
$ws.Range("F47").Value = 0.8961285909712722

$ws.Range("A48").Value = "7"
$ws.Range("B48").Value = "SVC"
$ws.Range("C48").Value = "Branch and Bound - Sem Normalização - Balanceado"
$ws.Range("D48").Value = "('sigmoid', 1, '0.0201')"
$ws.Range("E48").Value = "Recall"
$ws.Range("F48").Value = 0.02

$ws.Range("A49").Value = "7"
$ws.Range("B49").Value = "SVC"
$ws.Range("C49").Value = "Branch and Bound - Sem Normalização - Balanceado"
$ws.Range("D49").Value = "('sigmoid', 1, '0.0201')"
$ws.Range("E49").Value = "Precision"
$ws.Range("F49").Value = 0.006666666666666668
